$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$colC = @(4.927909215762419,4.750711121090916,4.640096136786081,4.594643996562552,4.587076506565086,4.639484562621838,4.867234215169501,5.296259412424064,5.596955309575498,5.729934621946245,5.779695868351556,5.76900605608611,5.734040688071694,5.712544578162932,5.588183967418262,5.510881612523719,5.466063290532977,5.45082896668603,5.519147792434892,5.744327365249833,5.888008340429301,5.811656343170477,5.515411820516682,5.182497424841224)
$colD = @(2.847003690747064,2.842352512310008,2.83951826737273,2.83836839667419,2.838177760593211,2.839502739525875,2.845395285006093,2.857137096684656,2.865899395334169,2.86991985634056,2.871447636775247,2.87111835948621,2.87004543903816,2.869388947936034,2.865637373057485,2.863344984699761,2.862029683714377,2.861584887426504,2.863588677445582,2.870360434484887,2.874817346523376,2.872435633780922,2.863478495811163,2.85393841118836)
$colE = @(16.70553859340209,15.74098122386723,15.12306693296054,14.86507095140317,14.82186579043061,15.11961218647729,16.37843040157255,18.78811643531084,20.46965862704372,21.19208178177448,21.45956983581151,21.40223123932992,21.21420978851935,21.0982507692216,20.42159408957464,19.99562288820606,19.74661129124329,19.66161140627439,20.04138240899331,21.2696009113484,22.03690264483969,21.63060736680953,20.02070734994359,18.13092776440014)
$colF = @(21.44205980522251,20.8346241685243,20.45982776462157,20.30690014922027,20.28150283021154,20.45776576031804,21.23315037590907,22.72822437241975,23.79757136625619,24.27519565339527,24.45460272838351,24.41603192272994,24.28998584363595,24.21258366221598,23.76616327202328,23.48989760279737,23.33017627434369,23.27596249551335,23.51939291460233,24.32704965068439,24.8463207816222,24.57001917168681,23.5060608495497,22.32795688141727)
$colG = @(27.28524728152907,25.9987004640269,25.18411712988713,24.84649731098813,24.79011010989612,25.17958608262031,26.84705228545231,29.90193584892394,31.99356861721881,32.90844354024771,33.24938747018263,33.1762069865081,32.93660461711536,32.78911880636145,31.93302164408605,31.39827571243353,31.08726844452471,30.98138443994495,31.45555771327269,33.00713256012471,33.9890400606857,33.46798535939458,31.42967163699104,29.10099330607954)
$colH = @(11.83963547182538,11.73936962845708,11.68066401629609,11.65748202894917,11.65367805690799,11.68034834722579,11.8044844931332,12.06954073893375,12.27596704518579,12.37207576034217,12.40875721936223,12.40084491739519,12.37508798279197,12.35934762366342,12.26972793351929,12.21529267921022,12.18419226319308,12.1736990455663,12.2210659575439,12.38264586451355,12.48990748496818,12.43251800643282,12.21845524834621,11.99567085756481)
$colO = @(18.21424911005318,17.82091217572675,17.58096086381734,17.48372304092086,17.46761413378029,17.57964707429929,18.07839353513229,19.06226355741023,19.78022120381482,20.10410721387772,20.22623451673262,20.19995739655341,20.11416586801241,20.06154445023205,19.75898847577349,19.57258892858084,19.46512847085713,19.42870566065304,19.59245815574105,20.13938006887658,20.49373176747096,20.30493171273612,19.58347619283926,18.7963700842242)

$colMap = @{ "C" = $colC; "D" = $colD; "E" = $colE; "F" = $colF; "G" = $colG; "H" = $colH; "O" = $colO }

foreach ($col in $colMap.Keys) {
    $values = $colMap[$col]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Range("$col$row").Value2 = $values[$i]
    }
}

